$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B4 from 2 to 0
$ws.Range("B2:B4").Value = 0

# Update C2:C15 timestamps from 13:26:35 to 13:26:47
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = "2025-04-04 13:26:47"
}
